$d = $word.ActiveDocument

$newText = "havainnointijaksot vuonna Hercules: 13.-22. Kesäkuuta, 12.-21. Heinäkuuta, 10.-19. Elokuuta"

# Find every paragraph whose text still carries the old "Perseus ... havainnointijaksot
# vuonna 2018 ..." date-range sentence (one occurrence also drags along a trailing "Ennen
# kuin menet ulos..." sentence + hyperlink). We collect the matching paragraph indices first
# and only then mutate, so that mid-loop reflows/Find state can't shift what "current" means.
$targets = New-Object System.Collections.ArrayList
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*havainnointijaksot vuonna 2018*") {
        [void]$targets.Add($i)
    }
}

foreach ($i in $targets) {
    $p = $d.Paragraphs.Item($i)
    $r = $p.Range
    # Exclude the trailing paragraph mark from the range so the paragraph itself survives.
    $r.End = $r.End - 1
    $r.Delete()
    $r.InsertAfter($newText)
}

Write-Host "Replaced $($targets.Count) paragraphs"
